$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: fill column A (Klasse) for the 21 new classes 06A..10D first.
# This matches the shared-string insertion order seen in the target file
# (all class codes were typed in one pass before the other columns).
$classesA = @("06A","06B","06C","06D","06E","07A","07B","07C","07D","08A","08B","08C","08D","09A","09B","09C","09D","10A","10B","10C","10D")
for ($i = 0; $i -lt $classesA.Count; $i++) {
    $r = 9 + $i
    $ws.Cells.Item($r, 1).Value = $classesA[$i]
}

# Step 2: fill the remaining columns (Sitzungsleiter, Von, Bis, Datum) for
# those same 21 rows (rows 9-29, before the later two-row insert).
$rowData = @{
    9  = @("Livia Schleßing, OStRin", "15:00", "15:30", "15.07.2017")
    10 = @("Livia Schleßing, OStRin", "15:00", "15:30", "15.07.2017")
    11 = @("Livia Schleßing, OStRin", "15:00", "15:30", "15.07.2017")
    12 = @("Martin Pabst, StD",       "14:00", "14:30", "15.07.2017")
    13 = @("Gerhard Maier, OStD",     "14:00", "14:30", "15.07.2017")
    14 = @("Martin Pabst, StD",       "14:00", "14:30", "15.07.2017")
    15 = @("Gerhard Maier, OStD",     "14:00", "14:30", "15.07.2017")
    16 = @("Martin Pabst, StD",       "14:00", "14:30", "15.07.2017")
    17 = @("Gerhard Maier, OStD",     "14:00", "14:30", "15.07.2017")
    18 = @("Andrea Fischer, StDin",   "14:30", "15:00", "15.07.2017")
    19 = @("Andrea Fischer, StDin",   "14:30", "15:00", "15.07.2017")
    20 = @("Andrea Fischer, StDin",   "14:30", "15:00", "15.07.2017")
    21 = @("Andrea Fischer, StDin",   "14:30", "15:00", "15.07.2017")
    22 = @("Andrea Fischer, StDin",   "14:30", "15:00", "15.07.2017")
    23 = @("Livia Schleßing, OStRin", "15:00", "15:30", "15.07.2017")
    24 = @("Livia Schleßing, OStRin", "15:00", "15:30", "15.07.2017")
    25 = @("Livia Schleßing, OStRin", "15:00", "15:30", "15.07.2017")
    26 = @("Martin Pabst, StD",       "14:00", "14:30", "15.07.2017")
    27 = @("Gerhard Maier, OStD",     "14:00", "14:30", "15.07.2017")
    28 = @("Martin Pabst, StD",       "14:00", "14:30", "15.07.2017")
    29 = @("Gerhard Maier, OStD",     "14:00", "14:30", "15.07.2017")
}

foreach ($r in 9..29) {
    $vals = $rowData[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
}

# Step 3: insert two new rows for classes 07E and 07F between 07D (row 17)
# and 08A (which is currently row 18, and will shift to row 20).
$ws.Rows("18:19").Insert()

$ws.Cells.Item(18, 1).Value = "07E"
$ws.Cells.Item(18, 2).Value = "Martin Pabst, StD"
$ws.Cells.Item(18, 3).Value = "14:00"
$ws.Cells.Item(18, 4).Value = "14:30"
$ws.Cells.Item(18, 5).Value = "15.07.2017"

$ws.Cells.Item(19, 1).Value = "07F"
$ws.Cells.Item(19, 2).Value = "Gerhard Maier, OStD"
$ws.Cells.Item(19, 3).Value = "14:00"
$ws.Cells.Item(19, 4).Value = "14:30"
$ws.Cells.Item(19, 5).Value = "15.07.2017"

# Step 4: update the selection/active cell to match the saved view state.
$ws.Range("A14").Select()
